# Scheduled market-data refresh: update the cached Universalis price/profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) on each class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4921.795
$ws.Range("I29").Value = 1950
$ws.Range("K29").Value = 5850
$ws.Range("M29").Value = -5569

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()

$ws.Range("H92").Value = 368.41666
$ws.Range("I92").Value = 231.5
$ws.Range("J92").Value = 642.25
$ws.Range("K92").Value = 231.5
$ws.Range("L92").Value = 642.25
$ws.Range("M92").Value = 1016.5
$ws.Range("N92").Value = -3138.25

$ws.Range("H96").Value = 920.7273
$ws.Range("J96").Value = 878
$ws.Range("L96").Value = 2634
$ws.Range("N96").Value = -5380

$ws.Range("H115").Value = 831.9286
$ws.Range("I115").Value = 742.0769
$ws.Range("K115").Value = 2226.2307
$ws.Range("M115").Value = -659.2307000000001

$ws.Range("H116").Value = 24298.4
$ws.Range("I116").Value = 27743
$ws.Range("J116").Value = 22002
$ws.Range("K116").Value = 27743
$ws.Range("L116").Value = 22002
$ws.Range("M116").Value = -24301
$ws.Range("N116").Value = -28886

$ws.Range("H132").Value = 3139.1292
$ws.Range("I132").Value = 2689.75
$ws.Range("K132").Value = 8069.25
$ws.Range("M132").Value = -5539.25

$ws.Range("H138").Value = 6600.547
$ws.Range("J138").Value = 6724.909
$ws.Range("L138").Value = 20174.727
$ws.Range("N138").Value = -30454.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H70").Value = 49998.5
$ws.Range("J70").Value = 49998.5
$ws.Range("L70").Value = 49998.5
$ws.Range("N70").Value = -50538.5

$ws.Range("H73").Value = 49998.5
$ws.Range("J73").Value = 49998.5
$ws.Range("L73").Value = 49998.5
$ws.Range("N73").Value = -51870.5

$ws.Range("H74").Value = 3293.1843
$ws.Range("I74").Value = 3069.1292
$ws.Range("J74").Value = 4285.4287
$ws.Range("K74").Value = 3069.1292
$ws.Range("L74").Value = 4285.4287
$ws.Range("M74").Value = -2195.1292
$ws.Range("N74").Value = -6033.4287

$ws.Range("H77").Value = 3293.1843
$ws.Range("I77").Value = 3069.1292
$ws.Range("J77").Value = 4285.4287
$ws.Range("K77").Value = 15345.646
$ws.Range("L77").Value = 21427.1435
$ws.Range("M77").Value = -10977.646
$ws.Range("N77").Value = -30163.1435

$ws.Range("H132").Value = 5984
$ws.Range("I132").Value = 5313.7407
$ws.Range("K132").Value = 15941.2221
$ws.Range("M132").Value = -13411.2221

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 79999
$ws.Range("J57").Value = 79999
$ws.Range("L57").Value = 79999
$ws.Range("N57").Value = -81439

$ws.Range("H99").Value = 2428.8333
$ws.Range("J99").Value = 2024.75
$ws.Range("L99").Value = 2024.75
$ws.Range("N99").Value = -5020.75

$ws.Range("H101").Value = 36620
$ws.Range("J101").Value = 36620
$ws.Range("L101").Value = 36620
$ws.Range("N101").Value = -43110

$ws.Range("H102").Value = 6639
$ws.Range("I102").Value = 6639
$ws.Range("K102").Value = 6639
$ws.Range("M102").Value = -3394

$ws.Range("H105").Value = 3108.9167
$ws.Range("J105").Value = 2474.5
$ws.Range("L105").Value = 2474.5
$ws.Range("N105").Value = -5968.5

$ws.Range("H134").Value = 19233.371
$ws.Range("I134").Value = 3018.3208
$ws.Range("K134").Value = 9054.9624
$ws.Range("M134").Value = -6519.9624

$ws.Range("H136").Value = 79999
$ws.Range("J136").Value = 79999
$ws.Range("L136").Value = 79999
$ws.Range("N136").Value = -90199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31933.8
$ws.Range("I31").Value = 1046.4445
$ws.Range("J31").Value = 64638.06
$ws.Range("K31").Value = 1046.4445
$ws.Range("L31").Value = 64638.06
$ws.Range("M31").Value = -751.4445000000001
$ws.Range("N31").Value = -65228.06

$ws.Range("H34").Value = 31933.8
$ws.Range("I34").Value = 1046.4445
$ws.Range("J34").Value = 64638.06
$ws.Range("K34").Value = 1046.4445
$ws.Range("L34").Value = 64638.06
$ws.Range("M34").Value = -844.4445000000001
$ws.Range("N34").Value = -65042.06

$ws.Range("H93").Value = 9960.571
$ws.Range("I93").Value = 8287.333000000001
$ws.Range("J93").Value = 20000
$ws.Range("K93").Value = 8287.333000000001
$ws.Range("L93").Value = 20000
$ws.Range("M93").Value = -6415.333000000001
$ws.Range("N93").Value = -23744

$ws.Range("H99").Value = 72549.266
$ws.Range("I99").Value = 5781.375
$ws.Range("J99").Value = 148855.42
$ws.Range("K99").Value = 5781.375
$ws.Range("L99").Value = 148855.42
$ws.Range("M99").Value = -4283.375
$ws.Range("N99").Value = -151851.42

$ws.Range("H126").Value = 72549.266
$ws.Range("I126").Value = 5781.375
$ws.Range("J126").Value = 148855.42
$ws.Range("K126").Value = 17344.125
$ws.Range("L126").Value = 446566.26
$ws.Range("M126").Value = -14874.125
$ws.Range("N126").Value = -451506.26

$ws.Range("H141").Value = 214689.36
$ws.Range("J141").Value = 214229.8
$ws.Range("L141").Value = 214229.8
$ws.Range("N141").Value = -224589.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1250
$ws.Range("J12").Value = 1250
$ws.Range("L12").Value = 3750
$ws.Range("N12").Value = -4096

$ws.Range("H37").Value = 100695.43
$ws.Range("J37").Value = 100695.43
$ws.Range("L37").Value = 302086.29
$ws.Range("N37").Value = -302310.29

$ws.Range("H56").Value = 7015.727
$ws.Range("I56").Value = 7015.727
$ws.Range("K56").Value = 7015.727
$ws.Range("M56").Value = -6485.727

$ws.Range("H87").Value = 28989.615
$ws.Range("I87").Value = 25266.428
$ws.Range("K87").Value = 75799.284
$ws.Range("M87").Value = -74551.284

$ws.Range("H90").Value = 28989.615
$ws.Range("I90").Value = 25266.428
$ws.Range("K90").Value = 227397.852
$ws.Range("M90").Value = -221157.852

$ws.Range("H92").Value = 1449.6428
$ws.Range("J92").Value = 1416.2222
$ws.Range("L92").Value = 4248.6666
$ws.Range("N92").Value = -6744.6666

$ws.Range("H93").Value = 5026
$ws.Range("I93").Value = 1080
$ws.Range("J93").Value = 6999
$ws.Range("K93").Value = 3240
$ws.Range("L93").Value = 20997
$ws.Range("M93").Value = -1368
$ws.Range("N93").Value = -24741

$ws.Range("H95").Value = 20000
$ws.Range("I95").Value = 20000
$ws.Range("K95").Value = 60000
$ws.Range("M95").Value = -57941

$ws.Range("H100").Value = 4764
$ws.Range("J100").Value = 4764
$ws.Range("L100").Value = 14292
$ws.Range("N100").Value = -15914

$ws.Range("H116").Value = 202006
$ws.Range("J116").Value = 3016
$ws.Range("L116").Value = 9048
$ws.Range("N116").Value = -15932

$ws.Range("H129").Value = 41668772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 414731.12
$ws.Range("I132").Value = 392220.78
$ws.Range("K132").Value = 1176662.34
$ws.Range("M132").Value = -1174132.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5804.077
$ws.Range("I7").Value = 4250.6
$ws.Range("K7").Value = 4250.6
$ws.Range("M7").Value = -4138.6

$ws.Range("H36").Value = 100750
$ws.Range("J36").Value = 100750
$ws.Range("L36").Value = 100750
$ws.Range("N36").Value = -101874

$ws.Range("H40").Value = 6686.2856
$ws.Range("I40").Value = 4952
$ws.Range("K40").Value = 4952
$ws.Range("M40").Value = -4816

$ws.Range("H55").Value = 461.66666
$ws.Range("I55").Value = 376.85715
$ws.Range("J55").Value = 580.4
$ws.Range("K55").Value = 376.85715
$ws.Range("L55").Value = 580.4
$ws.Range("M55").Value = -203.85715
$ws.Range("N55").Value = -926.4

$ws.Range("H61").Value = 6742.6294
$ws.Range("I61").Value = 7576.5654
$ws.Range("J61").Value = 1947.5
$ws.Range("K61").Value = 7576.5654
$ws.Range("L61").Value = 1947.5
$ws.Range("M61").Value = -7374.5654
$ws.Range("N61").Value = -2351.5

$ws.Range("H100").Value = 2800
$ws.Range("I100").Value = 2800
$ws.Range("K100").Value = 2800
$ws.Range("M100").Value = -2259

$ws.Range("H113").Value = 6742.6294
$ws.Range("I113").Value = 7576.5654
$ws.Range("J113").Value = 1947.5
$ws.Range("K113").Value = 7576.5654
$ws.Range("L113").Value = 1947.5
$ws.Range("M113").Value = -5406.5654
$ws.Range("N113").Value = -6287.5

$ws.Range("H126").Value = 5804.077
$ws.Range("I126").Value = 4250.6
$ws.Range("K126").Value = 12751.8
$ws.Range("M126").Value = -10281.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 621.8889
$ws.Range("I113").Value = 628.75
$ws.Range("J113").Value = 567
$ws.Range("K113").Value = 1886.25
$ws.Range("L113").Value = 1701
$ws.Range("M113").Value = 283.75
$ws.Range("N113").Value = -6041

$ws.Range("H126").Value = 1571.1
$ws.Range("I126").Value = 1537.5
$ws.Range("K126").Value = 4612.5
$ws.Range("M126").Value = -2142.5

$ws.Range("H132").Value = 54847.844
$ws.Range("I132").Value = 2339.389
$ws.Range("K132").Value = 7018.167
$ws.Range("M132").Value = -4488.167
